# Automatische test-sync: 2025-06-18 11:30:10
# Appends a new mail-log entry to the "Logs" sheet and updates the
# "Dashboard" pivot-style summary sheet + its chart accordingly.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new row to the Logs sheet (row 8) -----------------
$logs.Cells.Item(8, 1).Value = "Wat zijn jullie openingstijden?"
$logs.Cells.Item(8, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(8, 3).Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Cells.Item(8, 4).Value = "Informatieaanvraag"
$logs.Cells.Item(8, 5).Value = "Beste klant,`nBedankt voor je bericht. Onze openingstijden zijn maandag tot en met vrijdag van 09:00 tot 18:00 uur. Voor meer informatie of vragen zijn wij bereikbaar via e-mail of telefoon.`nMet vriendelijke groet,`n[Naam van het bedrijf]"
$logs.Cells.Item(8, 6).Value = "2025-06-18 11:00:13"
$logs.Cells.Item(8, 7).Value = "Ja"

# The multi-line "Antwoord" text would otherwise make the engine pin an
# explicit (auto-computed) row height with customHeight="1"; AutoFit()
# re-measures it and clears the custom/explicit height again, keeping the
# row definition clean (no ht/customHeight attribute), just like every
# other row on this sheet.
$logs.Rows.Item(8).AutoFit()

# --- 2. Append the new summary row to the Dashboard sheet (row 6) ----
$dash.Cells.Item(6, 1).Value = "Informatieaanvraag"
$dash.Cells.Item(6, 2).Value = 1

# --- 3. Extend the conditional formatting ranges on Logs to include
#        the newly added row.
foreach ($cf in $logs.Range("D2:D7").FormatConditions) {
    $cf.ModifyAppliesToRange($logs.Range("D2:D8"))
}
foreach ($cf in $logs.Range("G2:G7").FormatConditions) {
    $cf.ModifyAppliesToRange($logs.Range("G2:G8"))
}

# --- 4. Extend the chart series ranges on the Dashboard chart to
#        include the newly added category/value.
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$6"
$series.Values = "='Dashboard'!`$B`$2:`$B`$6"
